$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 99
$ws.Range("F3").Value = 135
$ws.Range("H3").Value = 3
$ws.Range("F4").Value = 106
$ws.Range("F5").Value = 73
$ws.Range("F6").Value = 99
$ws.Range("F7").Value = 267
$ws.Range("F8").Value = 150
$ws.Range("F9").Value = 85
$ws.Range("H9").Value = 8
$ws.Range("F10").Value = 80
$ws.Range("F11").Value = 104
$ws.Range("F12").Value = 319
$ws.Range("F13").Value = 302
$ws.Range("F14").Value = 151
$ws.Range("F15").Value = 265
$ws.Range("F16").Value = 316
$ws.Range("G16").Value = 1
$ws.Range("F17").Value = 285
$ws.Range("F18").Value = 310
$ws.Range("F19").Value = 510
$ws.Range("F20").Value = 301
$ws.Range("F21").Value = 330
$ws.Range("F22").Value = 1859
$ws.Range("F23").Value = 642
$ws.Range("F24").Value = 912
$ws.Range("F25").Value = 371
$ws.Range("F26").Value = 1934
$ws.Range("F27").Value = 506
$ws.Range("F28").Value = 764
$ws.Range("H28").Value = 8
$ws.Range("F29").Value = 1939
$ws.Range("F30").Value = 936
$ws.Range("F31").Value = 680
$ws.Range("F32").Value = 626
$ws.Range("G32").Value = 3
$ws.Range("F33").Value = 421
$ws.Range("F34").Value = 468
$ws.Range("F35").Value = 799
$ws.Range("F36").Value = 1847
$ws.Range("F37").Value = 549
$ws.Range("F38").Value = 659
$ws.Range("F39").Value = 459
$ws.Range("F40").Value = 1036
$ws.Range("F41").Value = 422
$ws.Range("F42").Value = 557
$ws.Range("H42").Value = 13
$ws.Range("F43").Value = 506
$ws.Range("G43").Value = 2
$ws.Range("F44").Value = 750
$ws.Range("F45").Value = 939
$ws.Range("F46").Value = 481
$ws.Range("F47").Value = 701
$ws.Range("G47").Value = 2
$ws.Range("F48").Value = 1832
$ws.Range("F49").Value = 2314
$ws.Range("F50").Value = 1672
$ws.Range("F51").Value = 4581
$ws.Range("F52").Value = 626
$ws.Range("F53").Value = 2525
$ws.Range("F54").Value = 768
$ws.Range("H54").Value = 61
$ws.Range("F55").Value = 3883
$ws.Range("G55").Value = 10
$ws.Range("F56").Value = 4024
$ws.Range("F57").Value = 1327
$ws.Range("F58").Value = 704
$ws.Range("G58").Value = 2
$ws.Range("F59").Value = 711
$ws.Range("F60").Value = 1529
$ws.Range("F61").Value = 6992
$ws.Range("G61").Value = 5
$ws.Range("F62").Value = 808
$ws.Range("F63").Value = 1055
$ws.Range("F64").Value = 1106
$ws.Range("F65").Value = 3950
$ws.Range("G65").Value = 22
$ws.Range("F66").Value = 1217
$ws.Range("F67").Value = 1304
$ws.Range("F68").Value = 1634
$ws.Range("F69").Value = 1230
$ws.Range("F70").Value = 1303
$ws.Range("F71").Value = 1455
$ws.Range("F72").Value = 1325
$ws.Range("F73").Value = 3423
$ws.Range("F74").Value = 4161
$ws.Range("F75").Value = 1783
$ws.Range("F76").Value = 1430
$ws.Range("F77").Value = 1551
$ws.Range("G77").Value = 6
$ws.Range("F78").Value = 3383
$ws.Range("F79").Value = 1309
$ws.Range("F80").Value = 8953
$ws.Range("F81").Value = 2087
$ws.Range("F82").Value = 1768
$ws.Range("F83").Value = 1373
$ws.Range("F84").Value = 5189
$ws.Range("F85").Value = 1532
$ws.Range("F86").Value = 1250
$ws.Range("F87").Value = 2874
$ws.Range("F88").Value = 3929
$ws.Range("F89").Value = 1614
$ws.Range("F90").Value = 3297
$ws.Range("G90").Value = 16
$ws.Range("F91").Value = 1978
$ws.Range("F92").Value = 1946
$ws.Range("G92").Value = 5
$ws.Range("F93").Value = 2107
$ws.Range("F94").Value = 2813
$ws.Range("F95").Value = 1593
$ws.Range("F96").Value = 2440
$ws.Range("G96").Value = 12
$ws.Range("F97").Value = 1914
$ws.Range("F98").Value = 1519
$ws.Range("G98").Value = 3
$ws.Range("F99").Value = 1894
$ws.Range("G99").Value = 8
$ws.Range("F100").Value = 4203
$ws.Range("F101").Value = 1300
$ws.Range("F102").Value = 2339
$ws.Range("G102").Value = 4
$ws.Range("F103").Value = 1492
$ws.Range("F104").Value = 6034
$ws.Range("G104").Value = 54
$ws.Range("F105").Value = 1399
$ws.Range("F106").Value = 1289
$ws.Range("F107").Value = 1298
$ws.Range("F108").Value = 2176
$ws.Range("F109").Value = 1918
$ws.Range("F110").Value = 1162
$ws.Range("F111").Value = 3564
$ws.Range("F112").Value = 1810
$ws.Range("H112").Value = 6
$ws.Range("F113").Value = 3185
$ws.Range("F114").Value = 2432
$ws.Range("F115").Value = 1420
$ws.Range("F116").Value = 3512
$ws.Range("G116").Value = 26
$ws.Range("F117").Value = 2604
$ws.Range("F118").Value = 2320
$ws.Range("F119").Value = 1559
$ws.Range("F120").Value = 2991
$ws.Range("G120").Value = 10
$ws.Range("F121").Value = 4477
$ws.Range("F122").Value = 5081
$ws.Range("F123").Value = 1539
$ws.Range("F124").Value = 2218
$ws.Range("F125").Value = 2361
$ws.Range("F126").Value = 2278
$ws.Range("F127").Value = 2257
$ws.Range("G127").Value = 19
$ws.Range("F128").Value = 2228
$ws.Range("F129").Value = 1918
$ws.Range("F130").Value = 4519
$ws.Range("F131").Value = 2175
$ws.Range("F132").Value = 2022
$ws.Range("F133").Value = 2530
$ws.Range("F134").Value = 5475
$ws.Range("F135").Value = 5145
$ws.Range("F136").Value = 4642
$ws.Range("F137").Value = 5282
$ws.Range("F138").Value = 4416
$ws.Range("G138").Value = 11
$ws.Range("F139").Value = 1661
$ws.Range("F140").Value = 2360
$ws.Range("F141").Value = 3853
$ws.Range("F142").Value = 5604
$ws.Range("F143").Value = 4617
$ws.Range("F144").Value = 5861
$ws.Range("G144").Value = 15
$ws.Range("F145").Value = 2964
$ws.Range("F146").Value = 7581
$ws.Range("G146").Value = 6
$ws.Range("F147").Value = 6301
$ws.Range("G147").Value = 12
